$d = $word.ActiveDocument

# Edit 1: ANOVA bullet point - rework wording about timing effects / aggregation / Tukey intervals
$d.Content.Find.Execute(
    "between the different years and months to consider the effects that timing has on sales. Revenue data would be aggregated by both month and year, and then the difference in means considered for each of these, with the Tukey intervals used to ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "between the different years and months to consider their effects on sales. Revenue data would be aggregated by month and year, and then the difference in means considered for each of these. The Tukey intervals would be used to ",
    2)

# Edit 2: remove "in this form of analysis" qualifier
$d.Content.Find.Execute(
    "(considered as factors in this form of analysis), with the final prediction model built around the p-values of these intervals.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(considered as factors), with the final prediction model built around the p-values of these intervals.",
    2)

# Edit 3: append justification sentence after "forecasts are found."
$d.Content.Find.Execute(
    "to test for goodness of fit before forecasts are found.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "to test for goodness of fit before forecasts are found. This model was used as it is likely to have the best trade-off between accuracy and time required as well as allowing for the best manual adjustment for the seasonal trends.",
    2)

# Edit 4: append closing sentence after the "Dynamic regression models..." paragraph
$d.Content.Find.Execute(
    "to see if more insights can be found. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "to see if more insights can be found. The other methods outlined above may also provide further insights.",
    2)
